$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "arreglo prestamo"
$ws.Range("D1").Value = "cuenta debito"

$ws.Range("E1").Select()
